$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.698.06'
$ws.Range('E2').Value = '  +0.42%  '
$ws.Range('D3').Value = '3.740.71'
$ws.Range('E3').Value = '  +0.70%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '612.38'
$ws.Range('E5').Value = '  +0.49%  '
$ws.Range('D6').Value = '178.58'
$ws.Range('E6').Value = '  +1.97%  '
$ws.Range('D7').Value = '3.741.94'
$ws.Range('E7').Value = '  +0.76%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -1.65%  '
$ws.Range('D10').Value = '0.165'
$ws.Range('E10').Value = '  +0.59%  '
$ws.Range('D11').Value = '6.60'
$ws.Range('E11').Value = '  +3.99%  '
$ws.Range('E12').Value = '  -2.93%  '
$ws.Range('D13').Value = '39.90'
$ws.Range('E13').Value = '  -1.27%  '
$ws.Range('D14').Value = '0.0000253'
$ws.Range('E14').Value = '  +0.51%  '
$ws.Range('D15').Value = '4.361.88'
$ws.Range('E15').Value = '  +0.83%  '
$ws.Range('D16').Value = '3.739.58'
$ws.Range('E16').Value = '  +0.71%  '
$ws.Range('D17').Value = '69.725.81'
$ws.Range('E17').Value = '  +0.40%  '
$ws.Range('E18').Value = '  -2.37%  '
$ws.Range('E19').Value = '  -1.12%  '
$ws.Range('D20').Value = '501.91'
$ws.Range('E20').Value = '  -1.83%  '
$ws.Range('D21').Value = '16.31'
$ws.Range('E21').Value = '  -2.15%  '
$ws.Range('E22').Value = '  -2.77%  '
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('D24').Value = '2.68'
$ws.Range('E24').Value = '  +9.36%  '
$ws.Range('D25').Value = '85.96'
$ws.Range('E25').Value = '  -1.66%  '
$ws.Range('D26').Value = '11.82'
$ws.Range('E26').Value = '  +8.38%  '
$ws.Range('D27').Value = '12.89'
$ws.Range('E27').Value = '  -3.16%  '
$ws.Range('D28').Value = '0.0000135'
$ws.Range('E28').Value = '  +8.75%  '
$ws.Range('E30').Value = '  -0.76%  '
$ws.Range('D33').Value = '30.36'
$ws.Range('E33').Value = '  -2.18%  '
$ws.Range('E34').Value = '  -0.91%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('E36').Value = '  +1.79%  '
$ws.Range('D37').Value = '6.11'
$ws.Range('E37').Value = '  -0.55%  '
$ws.Range('D38').Value = '0.357'
$ws.Range('E38').Value = '  +6.17%  '
$ws.Range('E39').Value = '  +4.11%  '
$ws.Range('D43').Value = '2.07'
$ws.Range('E43').Value = '  -4.05%  '
$ws.Range('E44').Value = '  -2.80%  '
$ws.Range('E45').Value = '  -2.27%  '
$ws.Range('D46').Value = '2.950.23'
$ws.Range('E46').Value = '  -4.05%  '
$ws.Range('D48').Value = '138.49'
$ws.Range('E48').Value = '  +2.80%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').Value = '27.13'
$ws.Range('E50').Value = '  -1.94%  '
$ws.Range('E51').Value = '  +0.40%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').Value = '8.12'
$ws.Range('E31').Value = '  +3.52%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '2.91'
$ws.Range('E32').Value = '  +3.31%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').Value = '450.44'
$ws.Range('E40').Value = '  +7.80%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').Value = '3.06'
$ws.Range('E41').Value = '  +13.63%  '
$ws.Range('B42').Value = 'Arweave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D42').Value = '46.30'
$ws.Range('E42').Value = '  +5.20%  '
